$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("threshold_b")
$ws.Range("A1").Select()
Write-Output "noop-ish"
